$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Expand the header row (row 1) from A:G to A:N -----------------------
# Existing B1:G1 already carry the bold/border header style (style 1);
# copy that format onto the newly added H1:N1 before filling their text so
# no new style entries are introduced.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Expand the data row (row 2) from A:G to A:N --------------------------
# B2:G2 already carry the plain data style (style 2); copy that format onto
# the newly added H2:N2 before filling their values.
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 1).Value = 33
$ws.Cells.Item(2, 2).Value = "HONDACRV2.4VTiS"
$ws.Cells.Item(2, 3).Value = 2354
$ws.Cells.Item(2, 4).Value = "潘瓊琪"
$ws.Cells.Item(2, 5).Value = "102年01月31曰"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = 950000
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"

# "date" (J2) must stay textual ("2013-12-31") instead of being
# auto-converted into a date serial, so force text format first.
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2013-12-31"

$ws.Cells.Item(2, 11).Value = "姚文智"
$ws.Cells.Item(2, 12).Value = 1745
$ws.Cells.Item(2, 13).Value = "tmpc2191"
$ws.Cells.Item(2, 14).Value = 33
